# Generate Report for Archive
# Update the handoff status text (shared by the Overview and per-locale sheets)
# from "Ready for handoff" to "In Translation", and shrink the status columns
# that were previously sized to fit the old (longer) status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update the status value everywhere it appears.
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# The status columns were sized for the old, longer text ("Ready for handoff");
# re-narrow them to fit the new, shorter text ("In Translation").
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
